# Commit: Fri, May 08, 2020 12:06:41 AM
#
# Update the table style applied to the three tables (on slides 14, 15
# and 16) from the deck's custom "Table_0" style to the built-in table
# style {90467875-FB9F-4FF7-8948-AB184F64548E}.

$p = $ppt.ActivePresentation
$newStyleId = "{90467875-FB9F-4FF7-8948-AB184F64548E}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
